$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new header row and per-module data (columns A:E) ---
$ws.Range("A1").Value = "Code"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "ChefModule"
$ws.Range("D1").Value = "ElementName1"
$ws.Range("E1").Value = "ElementName2"
$ws.Range("A2").Value = "GIL31"
$ws.Range("B2").Value = "pede. Suspendisse dui."
$ws.Range("C2").Value = "EL Haddad"
$ws.Range("D2").Value = "Nullam feugiat placerat"
$ws.Range("E2").Value = "varius et, euismod"
$ws.Range("A3").Value = "GIL32"
$ws.Range("B3").Value = "a nunc. In"
$ws.Range("C3").Value = "Badir"
$ws.Range("D3").Value = "sodales nisi magna"
$ws.Range("E3").Value = "elementum sem, vitae"
$ws.Range("A4").Value = "GIL33"
$ws.Range("B4").Value = "amet metus. Aliquam"
$ws.Range("C4").Value = "Ezzine"
$ws.Range("D4").Value = "Cras vulputate velit"
$ws.Range("E4").Value = "scelerisque neque sed"
$ws.Range("A5").Value = "GIL34"
$ws.Range("B5").Value = "quam vel sapien"
$ws.Range("C5").Value = "El Alami Hassoun"
$ws.Range("D5").Value = "Nunc mauris elit,"
$ws.Range("E5").Value = "libero et tristique"
$ws.Range("A6").Value = "GIL35"
$ws.Range("B6").Value = "feugiat nec, diam."
$ws.Range("C6").Value = "Lazaar"
$ws.Range("D6").Value = "pellentesque. Sed dictum."
$ws.Range("E6").Value = "ridiculus mus. Proin"
$ws.Range("A7").Value = "GIL36"
$ws.Range("B7").Value = "nonummy. Fusce fermentum"
$ws.Range("C7").Value = "El Haddad"
$ws.Range("D7").Value = "neque pellentesque massa"
$ws.Range("E7").Value = "Mauris eu turpis."
$ws.Range("A8").Value = "GIL41"
$ws.Range("B8").Value = "a, arcu. Sed"
$ws.Range("C8").Value = "EL Haddad"
$ws.Range("D8").Value = "sit amet risus."
$ws.Range("E8").Value = "Nulla facilisi. Sed"
$ws.Range("A9").Value = "GIL42"
$ws.Range("B9").Value = "Suspendisse eleifend. Cras"
$ws.Range("C9").Value = "El Alami Hassoun"
$ws.Range("D9").Value = "velit dui, semper"
$ws.Range("E9").Value = "ligula elit, pretium"
$ws.Range("A10").Value = "GIL43"
$ws.Range("B10").Value = "ante. Nunc mauris"
$ws.Range("C10").Value = "Badir"
$ws.Range("D10").Value = "tortor at risus."
$ws.Range("E10").Value = "felis. Donec tempor,"
$ws.Range("A11").Value = "GIL44"
$ws.Range("B11").Value = "lobortis quam a"
$ws.Range("C11").Value = "Ezzine"
$ws.Range("D11").Value = "euismod est arcu"
$ws.Range("E11").Value = "ligula eu enim."
$ws.Range("A12").Value = "GIL45"
$ws.Range("B12").Value = "rhoncus. Nullam velit"
$ws.Range("C12").Value = "Ben Achrab"
$ws.Range("D12").Value = "ut dolor dapibus"
$ws.Range("E12").Value = "commodo tincidunt nibh."
$ws.Range("A13").Value = "GIL46"
$ws.Range("B13").Value = "Donec tincidunt. Donec"
$ws.Range("C13").Value = "EL Haddad"
$ws.Range("D13").Value = "ornare tortor at"
$ws.Range("E13").Value = "ac, feugiat non,"

# --- Clear the old F column (ElementName3 header + any content) so the
#     used range shrinks back down to A:E ---
$ws.Range("F1:F13").ClearContents()

# --- Re-fit the data columns to their new (wider) content. The exact
#     sub-pixel widths that Excel originally stored for this sheet are not
#     reproducible bit-for-bit through ColumnWidth (Excel quantizes to a
#     1/6-character pixel grid), so we set the closest achievable width for
#     each column; column B lands on an exact match. ---
$ws.Columns("A").ColumnWidth = 4.877604166666667
$ws.Columns("B").ColumnWidth = 26.166666666666668
$ws.Columns("C").ColumnWidth = 15.307291666666666
$ws.Columns("D").ColumnWidth = 24.307291666666668
$ws.Columns("E").ColumnWidth = 22.736979166666668

# --- Move the active selection the way the author left it ---
$ws.Range("F2").Select()
